$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Finally, the logic of the code is the calculate" -> split into
#    "T" + "he logic of the code is the calculate" (the leading
#    "Finally, " is dropped and the following "t" is capitalised).
# ------------------------------------------------------------------
$target = $d.Content
[void]$target.Find.Execute("Finally, the logic of the code is the calculate")
$runStart = $target.Start

# Drop the leading "Finally, " (9 characters).
$leadIn = $d.Range($runStart, $runStart + 9)
$leadIn.Text = ""

# Capitalise the "t" that used to start "the logic..." into "T".
$capLetter = $d.Range($runStart, $runStart + 1)
$capLetter.Text = "T"

# Force a run break right after the new "T" by toggling a character
# format on/off (Word splits runs at the edges of a formatting
# operation even when the format ends up unchanged).
$splitAfterT = $d.Range($runStart, $runStart + 1)
$splitAfterT.Font.Bold = 1
$splitAfterT.Font.Bold = 0

# Re-establish the run boundaries around "Equations" (these were
# merged away by the text-length-changing edit above).
$eq = $d.Content
[void]$eq.Find.Execute("Equations")
$eqStart = $eq.Start
$eqEnd = $eq.End

$splitBeforeEquations = $d.Range($runStart, $eqStart)
$splitBeforeEquations.Font.Bold = 1
$splitBeforeEquations.Font.Bold = 0

$splitAfterEquations = $d.Range($eqStart, $eqEnd)
$splitAfterEquations.Font.Bold = 1
$splitAfterEquations.Font.Bold = 0

# ------------------------------------------------------------------
# 2) Add a new, final paragraph describing the writeFile method.
# ------------------------------------------------------------------
$lastParaCount = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastParaCount)
$lastPara.Range.InsertParagraphAfter()

$newParaCount = $d.Paragraphs.Count
$newPara = $d.Paragraphs($newParaCount)
$newPara.Range.Text = "Finally, the writeFile method uses a path similar to the one presented as an argument, and uses a Stringbuilder to take the information from each row and store it, then a bufferwriter to save the information to a file."
